# feat: add 2022-Q1 data
#
# 1. The worksheet that used to be named "总计" (sheetId=6) is renamed to
#    "2022-Q1" and its content is replaced with the fund-holdings detail
#    table for that quarter (same shape as the other quarterly sheets).
# 2. A brand new worksheet named "总计" (gets the next sheetId=7) is
#    inserted right after "2022-Q1" (i.e. at the end) and holds the
#    historical quarterly roll-up table, with a new row prepended for
#    "2022-Q1".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the new "2022-Q1" detail sheet
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# Use an existing quarterly sheet as a style template (bold/centered/
# bordered header row + bold/centered/bordered index column).
$template = $wb.Worksheets.Item("2021-Q4")

$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$template.Range("A2").Copy()
$q1.Range("A2:A27").PasteSpecial(-4122)

# Header row
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# Columns B-G are stored as text (even the numeric-looking ones, to
# preserve formatting such as leading zeros in fund codes and fixed
# decimal places); A and H are real numbers.
$q1.Range("B2:G27").NumberFormat = "@"

$rows = @(
    @("512800", "华宝中证银行ETF", "101.81", "99.20", "3.99", "4.0622", 9),
    @("510810", "汇添富中证上海国企ETF", "68.43", "99.71", "5.65", "3.8663", 6),
    @("515290", "天弘中证银行ETF", "84.75", "99.62", "4.02", "3.4070", 9),
    @("512700", "南方中证银行ETF", "25.07", "99.67", "4.01", "1.0053", 9),
    @("510230", "国泰上证180金融ETF", "36.29", "99.95", "2.62", "0.9508", 10),
    @("161121", "易方达中证银行指数（LOF）A", "20.63", "94.30", "3.79", "0.7819", 9),
    @("161723", "招商中证银行指数（LOF）", "18.28", "95.02", "3.81", "0.6965", 9),
    @("167301", "方正富邦中证保险主题指数（LOF）", "54.05", "93.05", "1.09", "0.5891", 9),
    @("161029", "富国中证银行指数", "15.35", "93.96", "3.78", "0.5802", 9),
    @("160631", "鹏华中证银行指数（LOF）", "12.80", "93.43", "3.76", "0.4813", 9),
    @("160517", "博时中证银行指数（LOF）", "10.14", "94.84", "3.81", "0.3863", 9),
    @("512820", "汇添富中证银行ETF", "6.97", "99.38", "4.01", "0.2795", 9),
    @("515020", "华夏中证银行ETF", "6.63", "99.17", "3.98", "0.2639", 9),
    @("009860", "易方达中证银行指数（LOF）C", "4.69", "94.30", "3.79", "0.1778", 9),
    @("159887", "富国中证800银行交易型开放式指数证券投资基金", "2.78", "99.42", "4.09", "0.1137", 8),
    @("160418", "华安中证银行指数（LOF）A", "2.86", "94.36", "3.80", "0.1087", 9),
    @("512730", "鹏华中证银行ETF", "1.99", "97.57", "3.91", "0.0778", 9),
    @("516310", "易方达中证银行交易型开放式指数证券投资基金", "1.82", "99.14", "3.99", "0.0726", 9),
    @("002849", "金信智能中国2025灵活配置混合", "1.44", "82.44", "4.50", "0.0648", 8),
    @("510760", "国泰上证综合ETF", "2.22", "95.43", "1.42", "0.0315", 8),
    @("011971", "西藏东财中证银行指数型发起式证券投资基金A", "0.80", "94.81", "3.81", "0.0305", 9),
    @("515280", "富国中证银行ETF", "0.67", "99.70", "4.02", "0.0269", 9),
    @("510650", "华夏金融ETF", "0.73", "99.00", "3.01", "0.0220", 10),
    @("515500", "海富通中证长三角领先ETF", "0.42", "95.16", "4.71", "0.0198", 5),
    @("011972", "西藏东财中证银行指数型发起式证券投资基金C", "0.35", "94.81", "3.81", "0.0133", 9),
    @("168205", "中融中证银行指数（LOF）", "0.35", "92.51", "3.73", "0.0131", 9)
)

$r = 2
foreach ($row in $rows) {
    $q1.Cells.Item($r, 1).Value = $r - 2
    $q1.Cells.Item($r, 2).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = $row[2]
    $q1.Cells.Item($r, 5).Value = $row[3]
    $q1.Cells.Item($r, 6).Value = $row[4]
    $q1.Cells.Item($r, 7).Value = $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: insert the brand new "总计" roll-up sheet after "2022-Q1"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$template.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$template.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Cells.Item(1,2).Value = "日期"
$total.Cells.Item(1,3).Value = "持有数量(只)"
$total.Cells.Item(1,4).Value = "持有市值(亿元)"

$summaryRows = @(
    @("2022-Q1", 26, 18.12),
    @("2021-Q4", 29, 19.62),
    @("2021-Q3", 44, 24.23),
    @("2021-Q2", 40, 22.11),
    @("2021-Q1", 67, 25.68),
    @("2020-Q4", 39, 19.24)
)

$r = 2
foreach ($row in $summaryRows) {
    $total.Cells.Item($r, 1).Value = $r - 2
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}
